$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new measurement row (row 20) by duplicating the last existing row
# (row 19) - this carries down the same number formats/styles - then
# overwrite the copied values with the new reading's data, same as a user
# would do when logging a new weigh-in.
$ws.Range("A19:H19").Copy() | Out-Null
$ws.Range("A20:H20").Insert() | Out-Null

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 43226
$ws.Range("C20").Value = 14.12
$ws.Range("D20").Value = 94.3
$ws.Range("E20").Value = 208
$ws.Range("H20").Value = 20.7

# F20/G20 already carry the copied formulas (E20-E19 / ROUND((D20/1.88)/1.88,2))
# from the Insert above, extending the shared formulas down one row.

# Update the view: select the next empty data row and scroll back to the top.
$ws.Range("A1").Select() | Out-Null
$ws.Range("F23").Select() | Out-Null
